# Auto-generated edit script applying cryptos price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.805.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.624.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -1.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.623.90"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.135"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.72%  "
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("E15").Value = "  +4.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.101.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.757.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.628.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "370.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.95%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("E23").Value = "  -1.93%  "
$ws.Range("E24").Value = "  -1.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.766.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E29").Value = "  +2.19%  "
$ws.Range("E30").Value = "  -3.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "576.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.18%  "
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "159.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("E40").Value = "  +4.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.368"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E42").Value = "  +1.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₆0335"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +15.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.81%  "
$ws.Range("E45").Value = "  +6.20%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "155.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.11%  "
$ws.Range("E51").Value = "  -0.65%  "
